$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (prices) that must stay text,
# exactly like the source inlineStr cells. Forcing NumberFormat to "@"
# before the write stops Excel from auto-converting them to numbers;
# ClearFormats() afterwards drops the temporary format so the cell keeps
# its original (default) style, same as every other cell on the sheet.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '27.158.80'
$ws.Range('E2').Value = '  +0.36%  '

Set-TextValue $ws.Range('D3') '1.676.66'
$ws.Range('E3').Value = '  -0.23%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue $ws.Range('D5') '214.62'
$ws.Range('E5').Value = '  -0.60%  '

Set-TextValue $ws.Range('D6') '0.518'
$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('E7').Value = '  +0.08%  '

Set-TextValue $ws.Range('D8') '22.90'
$ws.Range('E8').Value = '  +7.02%  '

Set-TextValue $ws.Range('D9') '0.261'
$ws.Range('E9').Value = '  +2.96%  '

Set-TextValue $ws.Range('D10') '0.0621'
$ws.Range('E10').Value = '  -0.51%  '

Set-TextValue $ws.Range('D11') '0.0889'
$ws.Range('E11').Value = '  +0.07%  '

Set-TextValue $ws.Range('D12') '1.914.52'
$ws.Range('E12').Value = '  -0.16%  '

Set-TextValue $ws.Range('D13') '1.674.64'
$ws.Range('E13').Value = '  -0.45%  '

Set-TextValue $ws.Range('D14') '4.21'
$ws.Range('E14').Value = '  +2.39%  '

Set-TextValue $ws.Range('D15') '0.562'
$ws.Range('E15').Value = '  +5.27%  '

Set-TextValue $ws.Range('D16') '66.58'
$ws.Range('E16').Value = '  +0.39%  '

Set-TextValue $ws.Range('D17') '27.129.67'
$ws.Range('E17').Value = '  +0.31%  '

Set-TextValue $ws.Range('D18') '235.47'
$ws.Range('E18').Value = '  -0.30%  '

Set-TextValue $ws.Range('D19') '7.91'
$ws.Range('E19').Value = '  -3.09%  '

Set-TextValue $ws.Range('D20') '0.0₃0741'
$ws.Range('E20').Value = '  +0.47%  '

$ws.Range('E21').Value = '  +0.00%  '

Set-TextValue $ws.Range('D22') '4.55'
$ws.Range('E22').Value = '  +1.71%  '

Set-TextValue $ws.Range('D23') '9.55'
$ws.Range('E23').Value = '  +3.06%  '

$ws.Range('E24').Value = '  -1.99%  '

Set-TextValue $ws.Range('D25') '148.47'
$ws.Range('E25').Value = '  +0.94%  '

Set-TextValue $ws.Range('D26') '7.48'
$ws.Range('E26').Value = '  +2.82%  '

Set-TextValue $ws.Range('D27') '16.39'
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('E28').Value = '  -0.31%  '

$ws.Range('E29').Value = '  +0.02%  '

Set-TextValue $ws.Range('D30') '0.0501'
$ws.Range('E30').Value = '  +0.72%  '

$ws.Range('E31').Value = '  -0.38%  '

Set-TextValue $ws.Range('D32') '3.36'
$ws.Range('E32').Value = '  -0.19%  '

Set-TextValue $ws.Range('D33') '1.541.60'
$ws.Range('E33').Value = '  -0.20%  '

Set-TextValue $ws.Range('D34') '3.23'
$ws.Range('E34').Value = '  +1.19%  '

Set-TextValue $ws.Range('D35') '1.66'
$ws.Range('E35').Value = '  -3.57%  '

Set-TextValue $ws.Range('D36') '0.609'
$ws.Range('E36').Value = '  +3.21%  '

Set-TextValue $ws.Range('D37') '0.945'
$ws.Range('E37').Value = '  +3.14%  '

$ws.Range('E38').Value = '  -0.05%  '

$ws.Range('E39').Value = '  -0.89%  '

$ws.Range('E40').Value = '  +2.02%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D41') '69.77'
$ws.Range('E41').Value = '  +2.79%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '5.78'
$ws.Range('E42').Value = '  +4.38%  '

$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('E44').Value = '  -0.39%  '

Set-TextValue $ws.Range('D45') '1.823.39'
$ws.Range('E45').Value = '  +0.08%  '

Set-TextValue $ws.Range('D46') '0.781'
$ws.Range('E46').Value = '  +0.27%  '

Set-TextValue $ws.Range('D47') '89.69'
$ws.Range('E47').Value = '  -0.93%  '

Set-TextValue $ws.Range('D48') '1.64'
$ws.Range('E48').Value = '  +6.71%  '

Set-TextValue $ws.Range('D49') '0.0₆0109'
$ws.Range('E49').Value = '  +1.42%  '

Set-TextValue $ws.Range('D50') '8.21'
$ws.Range('E50').Value = '  +2.84%  '

Set-TextValue $ws.Range('D51') '0.104'
$ws.Range('E51').Value = '  -0.21%  '
